$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"
